# Auto-generated: update currentAveragePrice / LevePrice / LeveProfit columns (H-N)
# across multiple sheets, per scheduled market-data refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
  # Row 6
  $ws.Range("H6").Value = 169.9
  $ws.Range("I6").Value = 157.16667
  $ws.Range("J6").Value = 284.5
  $ws.Range("K6").Value = 471.50001
  $ws.Range("L6").Value = 853.5
  $ws.Range("M6").Value = -359.50001
  $ws.Range("N6").Value = -1077.5
  # Row 17
  $ws.Range("H17").Value = 1840.4166
  $ws.Range("J17").Value = 1935
  $ws.Range("L17").Value = 5805
  $ws.Range("N17").Value = -6141
  # Row 18
  $ws.Range("H18").Value = 391.66666
  $ws.Range("I18").Value = 383.33334
  $ws.Range("J18").Value = 400
  $ws.Range("K18").Value = 383.33334
  $ws.Range("L18").Value = 400
  $ws.Range("M18").Value = -99.33334000000002
  $ws.Range("N18").Value = -968
  # Row 19
  $ws.Range("H19").Value = 244.8
  $ws.Range("I19").Value = 156.33333
  $ws.Range("J19").Value = 282.7143
  $ws.Range("K19").Value = 156.33333
  $ws.Range("L19").Value = 282.7143
  $ws.Range("M19").Value = 18.66667000000001
  $ws.Range("N19").Value = -632.7143
  # Row 40
  $ws.Range("H40").Value = 4359.8213
  $ws.Range("J40").Value = 2293
  $ws.Range("L40").Value = 2293
  $ws.Range("N40").Value = -2643
  # Row 107
  $ws.Range("H107").Value = 1490.5454
  $ws.Range("I107").Value = 1490.5454
  $ws.Range("K107").Value = 1490.5454
  $ws.Range("M107").Value = 429.4546
  # Row 116
  $ws.Range("H116").Value = 2819.1875
  $ws.Range("I116").Value = 2434.6667
  $ws.Range("K116").Value = 2434.6667
  $ws.Range("M116").Value = 1007.3333

$ws = $wb.Worksheets.Item("ARM")
  # Row 5
  $ws.Range("H5").Value = 268.1875
  $ws.Range("I5").Value = 268.23077
  $ws.Range("J5").Value = 268
  $ws.Range("K5").Value = 268.23077
  $ws.Range("L5").Value = 268
  $ws.Range("M5").Value = -156.23077
  $ws.Range("N5").Value = -492
  # Row 22
  $ws.Range("H22").Value = 1360.9
  $ws.Range("J22").Value = 1979.5
  $ws.Range("L22").Value = 1979.5
  $ws.Range("N22").Value = -2577.5
  # Row 101
  $ws.Range("H101").Value = 0
  $ws.Range("J101").Value = 0
  $ws.Range("L101").Value = 0
  $ws.Range("N101").Value = $null
  # Row 132
  $ws.Range("H132").Value = 6405.4
  $ws.Range("I132").Value = 6405.4
  $ws.Range("K132").Value = 19216.2
  $ws.Range("M132").Value = -16686.2

$ws = $wb.Worksheets.Item("BSM")
  # Row 4
  $ws.Range("H4").Value = 268.1875
  $ws.Range("I4").Value = 268.23077
  $ws.Range("J4").Value = 268
  $ws.Range("K4").Value = 268.23077
  $ws.Range("L4").Value = 268
  $ws.Range("M4").Value = -153.23077
  $ws.Range("N4").Value = -498
  # Row 20
  $ws.Range("H20").Value = 4278.8213
  $ws.Range("I20").Value = 5805.1763
  $ws.Range("K20").Value = 5805.1763
  $ws.Range("M20").Value = -5558.1763
  # Row 134
  $ws.Range("H134").Value = 4404.8335
  $ws.Range("I134").Value = 4404.8335
  $ws.Range("K134").Value = 13214.5005
  $ws.Range("M134").Value = -10679.5005

$ws = $wb.Worksheets.Item("CRP")
  # Row 16
  $ws.Range("H16").Value = 822.4
  $ws.Range("I16").Value = 734.7692
  $ws.Range("K16").Value = 734.7692
  $ws.Range("M16").Value = -447.7692
  # Row 60
  $ws.Range("H60").Value = 4990.6665
  $ws.Range("I60").Value = 4990.6665
  $ws.Range("K60").Value = 4990.6665
  $ws.Range("M60").Value = -4479.6665
  # Row 93
  $ws.Range("H93").Value = 21754.834
  $ws.Range("I93").Value = 11111
  $ws.Range("J93").Value = 23883.6
  $ws.Range("K93").Value = 11111
  $ws.Range("L93").Value = 23883.6
  $ws.Range("M93").Value = -9239
  $ws.Range("N93").Value = -27627.6
  # Row 103
  $ws.Range("H103").Value = 26888.572
  $ws.Range("I103").Value = 26370
  $ws.Range("K103").Value = 26370
  $ws.Range("M103").Value = -25198
  # Row 113
  $ws.Range("H113").Value = 822.4
  $ws.Range("I113").Value = 734.7692
  $ws.Range("K113").Value = 734.7692
  $ws.Range("M113").Value = 1435.2308
  # Row 132
  $ws.Range("H132").Value = 1526.5714
  $ws.Range("I132").Value = 1100
  $ws.Range("K132").Value = 3300
  $ws.Range("M132").Value = -770

$ws = $wb.Worksheets.Item("CUL")
  # Row 3
  $ws.Range("H3").Value = 0
  $ws.Range("I3").Value = 0
  $ws.Range("K3").Value = 0
  $ws.Range("M3").Value = $null
  # Row 28
  $ws.Range("H28").Value = 1416
  $ws.Range("I28").Value = 1624
  $ws.Range("J28").Value = 1000
  $ws.Range("K28").Value = 4872
  $ws.Range("L28").Value = 3000
  $ws.Range("M28").Value = -4640
  $ws.Range("N28").Value = -3464
  # Row 68
  $ws.Range("H68").Value = 1001
  $ws.Range("I68").Value = 0
  $ws.Range("J68").Value = 1001
  $ws.Range("K68").Value = 0
  $ws.Range("L68").Value = 3003
  $ws.Range("N68").Value = -4625
  $ws.Range("M68").Value = $null
  # Row 71
  $ws.Range("H71").Value = 1001
  $ws.Range("I71").Value = 0
  $ws.Range("J71").Value = 1001
  $ws.Range("K71").Value = 0
  $ws.Range("L71").Value = 9009
  $ws.Range("N71").Value = -17121
  $ws.Range("M71").Value = $null

$ws = $wb.Worksheets.Item("GSM")
  # Row 58
  $ws.Range("H58").Value = 9917.444
  $ws.Range("I58").Value = 6989
  $ws.Range("J58").Value = 10283.5
  $ws.Range("K58").Value = 6989
  $ws.Range("L58").Value = 10283.5
  $ws.Range("M58").Value = -6712
  $ws.Range("N58").Value = -10837.5
  # Row 103
  $ws.Range("H103").Value = 55000
  $ws.Range("J103").Value = 55000
  $ws.Range("L103").Value = 55000
  $ws.Range("N103").Value = -57344
  # Row 113
  $ws.Range("H113").Value = 3199.6667
  $ws.Range("J113").Value = 3283
  $ws.Range("L113").Value = 3283
  $ws.Range("N113").Value = -7623
  # Row 122
  $ws.Range("H122").Value = 3372.2778
  $ws.Range("I122").Value = 3888.818
  $ws.Range("J122").Value = 2560.5715
  $ws.Range("K122").Value = 11666.454
  $ws.Range("L122").Value = 7681.7145
  $ws.Range("M122").Value = -9216.454000000002
  $ws.Range("N122").Value = -12581.7145
  # Row 132
  $ws.Range("H132").Value = 3181
  $ws.Range("I132").Value = 3816.3333
  $ws.Range("J132").Value = 2799.8
  $ws.Range("K132").Value = 11448.9999
  $ws.Range("L132").Value = 8399.400000000001
  $ws.Range("M132").Value = -8918.999899999999
  $ws.Range("N132").Value = -13459.4

$ws = $wb.Worksheets.Item("LTW")
  # Row 40
  $ws.Range("H40").Value = 7141.909
  $ws.Range("I40").Value = 6119.2
  $ws.Range("K40").Value = 6119.2
  $ws.Range("M40").Value = -5983.2
  # Row 46
  $ws.Range("H46").Value = 3288.9412
  $ws.Range("I46").Value = 0
  $ws.Range("J46").Value = 3288.9412
  $ws.Range("K46").Value = 0
  $ws.Range("L46").Value = 3288.9412
  $ws.Range("N46").Value = -3664.9412
  $ws.Range("M46").Value = $null
  # Row 57
  $ws.Range("H57").Value = 8970.5
  $ws.Range("I57").Value = 8970.5
  $ws.Range("K57").Value = 8970.5
  $ws.Range("M57").Value = -8404.5
  # Row 68
  $ws.Range("H68").Value = 3555.4614
  $ws.Range("I68").Value = 3689.3
  $ws.Range("J68").Value = 3109.3333
  $ws.Range("K68").Value = 3689.3
  $ws.Range("L68").Value = 3109.3333
  $ws.Range("M68").Value = -2940.3
  $ws.Range("N68").Value = -4607.3333
  # Row 71
  $ws.Range("H71").Value = 3555.4614
  $ws.Range("I71").Value = 3689.3
  $ws.Range("J71").Value = 3109.3333
  $ws.Range("K71").Value = 18446.5
  $ws.Range("L71").Value = 15546.6665
  $ws.Range("M71").Value = -14702.5
  $ws.Range("N71").Value = -23034.6665
  # Row 82
  $ws.Range("H82").Value = 3035.0908
  $ws.Range("I82").Value = 2887.5
  $ws.Range("J82").Value = 3699.25
  $ws.Range("K82").Value = 2887.5
  $ws.Range("L82").Value = 3699.25
  $ws.Range("M82").Value = -2526.5
  $ws.Range("N82").Value = -4421.25
  # Row 85
  $ws.Range("H85").Value = 3035.0908
  $ws.Range("I85").Value = 2887.5
  $ws.Range("J85").Value = 3699.25
  $ws.Range("K85").Value = 2887.5
  $ws.Range("L85").Value = 3699.25
  $ws.Range("M85").Value = -1639.5
  $ws.Range("N85").Value = -6195.25
  # Row 93
  $ws.Range("H93").Value = 15162.827
  $ws.Range("I93").Value = 886
  $ws.Range("J93").Value = 60032.855
  $ws.Range("K93").Value = 886
  $ws.Range("L93").Value = 60032.855
  $ws.Range("M93").Value = 362
  $ws.Range("N93").Value = -62528.855

$ws = $wb.Worksheets.Item("WVR")
  # Row 56
  $ws.Range("H56").Value = 5000
  $ws.Range("I56").Value = 5000
  $ws.Range("K56").Value = 5000
  $ws.Range("M56").Value = -4286
  # Row 107
  $ws.Range("H107").Value = 2291
  $ws.Range("I107").Value = 1551.7142
  $ws.Range("J107").Value = 2937.875
  $ws.Range("K107").Value = 4655.142599999999
  $ws.Range("L107").Value = 8813.625
  $ws.Range("M107").Value = -2735.142599999999
  $ws.Range("N107").Value = -12653.625
  # Row 126
  $ws.Range("H126").Value = 2912.5
  $ws.Range("I126").Value = 2020.5555
  $ws.Range("K126").Value = 6061.666499999999
  $ws.Range("M126").Value = -3591.666499999999
  # Row 136
  $ws.Range("H136").Value = 726.15625
  $ws.Range("I136").Value = 637.9666999999999
  $ws.Range("K136").Value = 1913.9001
  $ws.Range("M136").Value = 636.0999000000002
